$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete extra rows 11-13 (target clusters reduced from 4 to 3 per sending cluster)
$ws.Range("A11:T13").EntireRow.Delete()

# Row 2: ECs -> FAPs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cntn2"
$ws.Cells.Item(2, 3).Value = "Nrcam"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.030141
$ws.Cells.Item(2, 8).Value = 0.090423
$ws.Cells.Item(2, 9).Value = 0.2387784203438168
$ws.Cells.Item(2, 10).Value = 0.2387784203438169
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1623456666666667
$ws.Cells.Item(2, 14).Value = 0.487037
$ws.Cells.Item(2, 15).Value = 0.07831009050961847
$ws.Cells.Item(2, 16).Value = 0.07831009050961847
$ws.Cells.Item(2, 17).Value = 0.004893260739
$ws.Cells.Item(2, 18).Value = 0.044039346651
$ws.Cells.Item(2, 19).Value = 0.01869875970886802
$ws.Cells.Item(2, 20).Value = 0.01869875970886803

# Row 3: ECs -> MuSCs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cntn2"
$ws.Cells.Item(3, 3).Value = "Nrcam"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.030141
$ws.Cells.Item(3, 8).Value = 0.090423
$ws.Cells.Item(3, 9).Value = 0.2387784203438168
$ws.Cells.Item(3, 10).Value = 0.2387784203438169
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.897446
$ws.Cells.Item(3, 14).Value = 5.692337999999999
$ws.Cells.Item(3, 15).Value = 0.9152641462380487
$ws.Cells.Item(3, 16).Value = 0.9152641462380487
$ws.Cells.Item(3, 17).Value = 0.057190919886
$ws.Cells.Item(3, 18).Value = 0.514718278974
$ws.Cells.Item(3, 19).Value = 0.2185453270360534
$ws.Cells.Item(3, 20).Value = 0.2185453270360535

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cntn2"
$ws.Cells.Item(4, 3).Value = "Nrcam"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.030141
$ws.Cells.Item(4, 8).Value = 0.090423
$ws.Cells.Item(4, 9).Value = 0.2387784203438168
$ws.Cells.Item(4, 10).Value = 0.2387784203438169
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.01332133333333333
$ws.Cells.Item(4, 14).Value = 0.039964
$ws.Cells.Item(4, 15).Value = 0.006425763252332764
$ws.Cells.Item(4, 16).Value = 0.006425763252332764
$ws.Cells.Item(4, 17).Value = 0.000401518308
$ws.Cells.Item(4, 18).Value = 0.003613664772
$ws.Cells.Item(4, 19).Value = 0.001534333598895364
$ws.Cells.Item(4, 20).Value = 0.001534333598895365

# Row 5: FAPs -> FAPs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cntn2"
$ws.Cells.Item(5, 3).Value = "Nrcam"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.027635
$ws.Cells.Item(5, 8).Value = 0.082905
$ws.Cells.Item(5, 9).Value = 0.2189257704190763
$ws.Cells.Item(5, 10).Value = 0.2189257704190763
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.1623456666666667
$ws.Cells.Item(5, 14).Value = 0.487037
$ws.Cells.Item(5, 15).Value = 0.07831009050961847
$ws.Cells.Item(5, 16).Value = 0.07831009050961847
$ws.Cells.Item(5, 17).Value = 0.004486422498333334
$ws.Cells.Item(5, 18).Value = 0.040377802485
$ws.Cells.Item(5, 19).Value = 0.01714409689640582
$ws.Cells.Item(5, 20).Value = 0.01714409689640582

# Row 6: FAPs -> MuSCs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Cntn2"
$ws.Cells.Item(6, 3).Value = "Nrcam"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.027635
$ws.Cells.Item(6, 8).Value = 0.082905
$ws.Cells.Item(6, 9).Value = 0.2189257704190763
$ws.Cells.Item(6, 10).Value = 0.2189257704190763
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.897446
$ws.Cells.Item(6, 14).Value = 5.692337999999999
$ws.Cells.Item(6, 15).Value = 0.9152641462380487
$ws.Cells.Item(6, 16).Value = 0.9152641462380487
$ws.Cells.Item(6, 17).Value = 0.05243592021
$ws.Cells.Item(6, 18).Value = 0.47192328189
$ws.Cells.Item(6, 19).Value = 0.2003749083521229
$ws.Cells.Item(6, 20).Value = 0.200374908352123

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Cntn2"
$ws.Cells.Item(7, 3).Value = "Nrcam"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.027635
$ws.Cells.Item(7, 8).Value = 0.082905
$ws.Cells.Item(7, 9).Value = 0.2189257704190763
$ws.Cells.Item(7, 10).Value = 0.2189257704190763
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.01332133333333333
$ws.Cells.Item(7, 14).Value = 0.039964
$ws.Cells.Item(7, 15).Value = 0.006425763252332764
$ws.Cells.Item(7, 16).Value = 0.006425763252332764
$ws.Cells.Item(7, 17).Value = 0.0003681350466666667
$ws.Cells.Item(7, 18).Value = 0.00331321542
$ws.Cells.Item(7, 19).Value = 0.00140676517054754
$ws.Cells.Item(7, 20).Value = 0.00140676517054754

# Row 8: MuSCs -> FAPs
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Cntn2"
$ws.Cells.Item(8, 3).Value = "Nrcam"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.068454
$ws.Cells.Item(8, 8).Value = 0.205362
$ws.Cells.Item(8, 9).Value = 0.5422958092371069
$ws.Cells.Item(8, 10).Value = 0.5422958092371069
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.1623456666666667
$ws.Cells.Item(8, 14).Value = 0.487037
$ws.Cells.Item(8, 15).Value = 0.07831009050961847
$ws.Cells.Item(8, 16).Value = 0.07831009050961847
$ws.Cells.Item(8, 17).Value = 0.011113210266
$ws.Cells.Item(8, 18).Value = 0.100018892394
$ws.Cells.Item(8, 19).Value = 0.04246723390434463
$ws.Cells.Item(8, 20).Value = 0.04246723390434463

# Row 9: MuSCs -> MuSCs
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Cntn2"
$ws.Cells.Item(9, 3).Value = "Nrcam"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.068454
$ws.Cells.Item(9, 8).Value = 0.205362
$ws.Cells.Item(9, 9).Value = 0.5422958092371069
$ws.Cells.Item(9, 10).Value = 0.5422958092371069
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.897446
$ws.Cells.Item(9, 14).Value = 5.692337999999999
$ws.Cells.Item(9, 15).Value = 0.9152641462380487
$ws.Cells.Item(9, 16).Value = 0.9152641462380487
$ws.Cells.Item(9, 17).Value = 0.129887768484
$ws.Cells.Item(9, 18).Value = 1.168989916356
$ws.Cells.Item(9, 19).Value = 0.4963439108498723
$ws.Cells.Item(9, 20).Value = 0.4963439108498723

# Row 10: MuSCs -> Resolving-Mac
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Cntn2"
$ws.Cells.Item(10, 3).Value = "Nrcam"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.068454
$ws.Cells.Item(10, 8).Value = 0.205362
$ws.Cells.Item(10, 9).Value = 0.5422958092371069
$ws.Cells.Item(10, 10).Value = 0.5422958092371069
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.01332133333333333
$ws.Cells.Item(10, 14).Value = 0.039964
$ws.Cells.Item(10, 15).Value = 0.006425763252332764
$ws.Cells.Item(10, 16).Value = 0.006425763252332764
$ws.Cells.Item(10, 17).Value = 0.0009118985519999999
$ws.Cells.Item(10, 18).Value = 0.008207086968
$ws.Cells.Item(10, 19).Value = 0.00348466448288986
$ws.Cells.Item(10, 20).Value = 0.00348466448288986

